{"js": "// Replace the 100 equation answers in the 20x5 practice-problem table\n// (row-major order) with the new values from the commit, keeping\n// everything else (formatting, the date paragraph, etc.) untouched.\nconst newValues = [\"97-53=44\", \"38+37=75\", \"78-61=17\", \"73-33=40\", \"8-2=6\", \"44+3=47\", \"30+26=56\", \"32+40=72\", \"55+11=66\", \"91-68=23\", \"32+44=76\", \"14+44=58\", \"27-9=18\", \"79-73=6\", \"99-32=67\", \"50+2=52\", \"31+43=74\", \"61-57=4\", \"17+71=88\", \"64-50=14\", \"89-47=42\", \"25+71=96\", \"72+13=85\", \"89-7=82\", \"18+33=51\", \"80-36=44\", \"70-65=5\", \"6+65=71\", \"87-28=59\", \"59-13=46\", \"77-59=18\", \"28-16=12\", \"17-16=1\", \"53+41=94\", \"60+34=94\", \"88-59=29\", \"98-64=34\", \"91-64=27\", \"39-32=7\", \"26+47=73\", \"34-28=6\", \"31+18=49\", \"4+77=81\", \"90-64=26\", \"44+48=92\", \"83+13=96\", \"65+7=72\", \"7+58=65\", \"97-82=15\", \"46-11=35\", \"47+38=85\", \"12+1=13\", \"7+43=50\", \"28+9=37\", \"43+29=72\", \"60+28=88\", \"11+32=43\", \"34+48=82\", \"22+69=91\", \"64+33=97\", \"28-8=20\", \"4+91=95\", \"66-15=51\", \"11+74=85\", \"88-11=77\", \"49+35=84\", \"95-37=58\", \"93-15=78\", \"2+78=80\", \"74-62=12\", \"40+53=93\", \"91-65=26\", \"31+35=66\", \"38-19=19\", \"61+32=93\", \"66-24=42\", \"20-0=20\", \"12+63=75\", \"29-6=23\", \"4+95=99\", \"70+16=86\", \"40-5=35\", \"81+12=93\", \"91-9=82\", \"41+18=59\", \"79-18=61\", \"41-33=8\", \"31-16=15\", \"60-23=37\", \"8+40=48\", \"54-41=13\", \"62-51=11\", \"12+13=25\", \"10+31=41\", \"55+5=60\", \"21+74=95\", \"94-16=78\", \"34-18=16\", \"40-22=18\", \"32+28=60\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst rowCount = rows.items.length;\nconst colCount = 5;\n\n// Rewrite just the text of each cell's own Range (NOT the cell body),\n// which replaces only the run's text and keeps the existing paragraph\n// formatting (w:pPr/w:jc) and run formatting (w:rPr/w:rFonts/w:sz) intact.\nlet i = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    cell.getRange().insertText(newValues[i], Word.InsertLocation.replace);\n    i++;\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the 100 equation answers in the 20x5 practice-problem table\n# (row-major order) with the new values from the commit, keeping\n# everything else (formatting, the date paragraph, etc.) untouched.\n$newValues = @(\"97-53=44\", \"38+37=75\", \"78-61=17\", \"73-33=40\", \"8-2=6\", \"44+3=47\", \"30+26=56\", \"32+40=72\", \"55+11=66\", \"91-68=23\", \"32+44=76\", \"14+44=58\", \"27-9=18\", \"79-73=6\", \"99-32=67\", \"50+2=52\", \"31+43=74\", \"61-57=4\", \"17+71=88\", \"64-50=14\", \"89-47=42\", \"25+71=96\", \"72+13=85\", \"89-7=82\", \"18+33=51\", \"80-36=44\", \"70-65=5\", \"6+65=71\", \"87-28=59\", \"59-13=46\", \"77-59=18\", \"28-16=12\", \"17-16=1\", \"53+41=94\", \"60+34=94\", \"88-59=29\", \"98-64=34\", \"91-64=27\", \"39-32=7\", \"26+47=73\", \"34-28=6\", \"31+18=49\", \"4+77=81\", \"90-64=26\", \"44+48=92\", \"83+13=96\", \"65+7=72\", \"7+58=65\", \"97-82=15\", \"46-11=35\", \"47+38=85\", \"12+1=13\", \"7+43=50\", \"28+9=37\", \"43+29=72\", \"60+28=88\", \"11+32=43\", \"34+48=82\", \"22+69=91\", \"64+33=97\", \"28-8=20\", \"4+91=95\", \"66-15=51\", \"11+74=85\", \"88-11=77\", \"49+35=84\", \"95-37=58\", \"93-15=78\", \"2+78=80\", \"74-62=12\", \"40+53=93\", \"91-65=26\", \"31+35=66\", \"38-19=19\", \"61+32=93\", \"66-24=42\", \"20-0=20\", \"12+63=75\", \"29-6=23\", \"4+95=99\", \"70+16=86\", \"40-5=35\", \"81+12=93\", \"91-9=82\", \"41+18=59\", \"79-18=61\", \"41-33=8\", \"31-16=15\", \"60-23=37\", \"8+40=48\", \"54-41=13\", \"62-51=11\", \"12+13=25\", \"10+31=41\", \"55+5=60\", \"21+74=95\", \"94-16=78\", \"34-18=16\", \"40-22=18\", \"32+28=60\")\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $t.Cell($r, $c)\n        # Assigning Range.Text (not including the end-of-cell mark) replaces\n        # only the run's text, leaving the paragraph/run formatting intact.\n        $cell.Range.Text = $newValues[$i]\n        $i++\n    }\n}\n"}
